# Rebuild the "ESTADO DE CUENTA" detail table (rows 16-49).
#
# Commit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The underlying query/database that feeds this report was refreshed: the
# detail rows are regrouped by worker (MAURICIO ALBERTO OTERO BULA first,
# then JOHON ALFRED CORONELL BOLIVAR) and, within each worker, ordered by
# "Periodo Mora" descending (most recent period first). A few "Valor Mora"
# amounts for the most recent periods were also recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: TipoDoc, NumDoc, Nombre, Periodo, ValorMora, SalarioBasico
$rows = @(
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2105",25850,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2104",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2103",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2102",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2101",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2012",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2011",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2010",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2009",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2008",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2007",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2006",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2005",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2004",32312,807803),
    @("CC","1069483129","MAURICIO ALBERTO OTERO BULA","2003",32312,807803),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2105",27200,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2104",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2103",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2102",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2101",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2012",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2011",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2010",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2009",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2008",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2007",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2006",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2005",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2004",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2003",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2002",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","2001",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","1912",34000,850000),
    @("CC","73594227","JOHON ALFRED CORONELL BOLIVAR","1911",34000,850000)
)

$firstRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $firstRow + $i
    $data = $rows[$i]

    $ws.Range("B$r").Value = $data[0]
    $ws.Range("C$r").Value = $data[1]
    $ws.Range("D$r").Value = $data[2]
    $ws.Range("E$r").Value = $data[3]
    $ws.Range("F$r").Value = $data[4]
    $ws.Range("G$r").Value = $data[5]
}
